# Edit script for flower8-register-map.xlsx
# Implements:
#  - D68: updated note text + row height 30
#  - Rows 93-100: new "COINC TRIGGER-> threshold chN" column C (bold, no color),
#    shared "lowest byte..." note in column D, "register address moved..." note in F93
#  - New rows 135-262: extended beam-threshold register block (0x80-0xFF)
#  - Final selection on D136

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registers")
$ws.Activate()

# xlPasteFormats
$xlPasteFormats = -4122

### 1. New register rows 135-262 (addresses 0x80-0xFF) -- structure + addresses first
# Row 135: plain continuation row (address 0x80), no extra content yet.
$ws.Range("A135").Value = 128
$ws.Range("B135").Formula = '="x" & DEC2HEX(A135,2)'

# Row 136: beam-threshold block header/description row -- gets a top border
# (reusing the existing "border=2" look already used elsewhere in the sheet)
# plus word-wrap on the trailing blank cell.
$ws.Range("A136").Value = 129
$ws.Range("B136").Formula = '="x" & DEC2HEX(A136,2)'
$ws.Range("C136").Value = "Beam 0 Thresholds"

$ws.Range("A136:E136").Borders.Item(8).LineStyle = 1
$ws.Range("A136:E136").Borders.Item(8).Weight = 2
$ws.Range("F136").Borders.Item(8).LineStyle = 1
$ws.Range("F136").Borders.Item(8).Weight = 2
$ws.Range("F136").WrapText = $true

# Rows 137-150: mechanical address/hex columns only.
for ($r = 137; $r -le 150; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 7
    $ws.Cells.Item($r, 2).Formula = '="x" & DEC2HEX(A' + $r + ',2)'
}

# Row 151: "Beam Threshold Block" label marks the start of the next beam's regs.
$ws.Range("C151").Value = "Beam Threshold Block"

# Rows 152-165: mechanical.
for ($r = 152; $r -le 165; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 7
    $ws.Cells.Item($r, 2).Formula = '="x" & DEC2HEX(A' + $r + ',2)'
}

# Row 166: closing divider row -- bottom border across A:F, with wrap on F.
$ws.Range("A166").Value = 159
$ws.Range("B166").Formula = '="x" & DEC2HEX(A166,2)'

$ws.Range("A166:E166").Borders.Item(9).LineStyle = 1
$ws.Range("A166:E166").Borders.Item(9).Weight = 2
$ws.Range("F166").Borders.Item(9).LineStyle = 1
$ws.Range("F166").Borders.Item(9).Weight = 2
$ws.Range("F166").WrapText = $true

# Rows 167-262: mechanical address/hex columns only (addresses 0xA0-0xFF).
for ($r = 167; $r -le 262; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 7
    $ws.Cells.Item($r, 2).Formula = '="x" & DEC2HEX(A' + $r + ',2)'
}

# Back to row 136 to fill in the "max value" note now that the block below it exists.
$ws.Range("E136").Value = "0xffffff"

### 2. Rows 93-100: rename beam-threshold regs to per-channel COINC TRIGGER thresholds
# Column C loses its color formatting, picking up the bold "section header" look
# already used by column A/E in this custom-formatted row block (style index 2).
$ws.Range("A93").Copy()
$ws.Range("C93:C100").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("C93").Value = "COINC TRIGGER-> threshold ch0"
$ws.Range("D93").Value = "lowest byte : nominal trigger thresh ; middle byte : servo trigger threshold"
$ws.Range("F93").Value = "register address moved in FLOWER8"

$servoNote = "lowest byte : nominal trigger thresh ; middle byte : servo trigger threshold"
$ws.Range("C94").Value = "COINC TRIGGER-> threshold ch1"
$ws.Range("D94").Value = $servoNote
$ws.Range("C95").Value = "COINC TRIGGER-> threshold ch2"
$ws.Range("D95").Value = $servoNote
$ws.Range("C96").Value = "COINC TRIGGER-> threshold ch3"
$ws.Range("D96").Value = $servoNote
$ws.Range("C97").Value = "COINC TRIGGER-> threshold ch4"
$ws.Range("D97").Value = $servoNote
$ws.Range("C98").Value = "COINC TRIGGER-> threshold ch5"
$ws.Range("D98").Value = $servoNote
$ws.Range("C99").Value = "COINC TRIGGER-> threshold ch6"
$ws.Range("D99").Value = $servoNote
$ws.Range("C100").Value = "COINC TRIGGER-> threshold ch7"
$ws.Range("D100").Value = $servoNote

### 3. Row 68: updated functionality note + taller row for wrapped text
$ws.Range("D68").Value = "LSB-> pps trig enable, middle byte LSB+1-> phased trig enable,,middle byte LSB-> coinc trig enable high byte LSB->ext trig enable"
$ws.Rows.Item(68).RowHeight = 30

### 4. Last edit of the session: the explanatory note beside the new beam block.
$ws.Range("D136").Value = "MS 12 bits = Servo Thresholds, LS 12 bits Trigger Thresholds"

### 5. Leave the selection where the author left it.
$ws.Range("D136").Select()
